$d = $word.ActiveDocument

# Helper: append a brand-new paragraph after $afterPara, return the new Paragraph object.
function New-ParaAfter($afterPara) {
    $afterPara.Range.InsertParagraphAfter() | Out-Null
    $doc = $word.ActiveDocument
    return $doc.Paragraphs.Item($doc.Paragraphs.Count)
}

# Helper: add one run of text at the end of a paragraph's range.
function Add-Run($para, [string]$text) {
    $r = $para.Range
    $r.Collapse(0)
    $r.InsertAfter($text)
}

# The current last paragraph in the document ("...STARL examples"), already
# styled ListParagraph / numId=2 -- remember its list template so later list
# paragraphs can continue the very same list instance instead of minting a
# new one.
$lastIndex = $d.Paragraphs.Count
$starlPara = $d.Paragraphs.Item($lastIndex)
$listTemplate = $starlPara.Range.ListFormat.ListTemplate

# 1) Title paragraph: "Monday 28/03 02:30pm"
$p = New-ParaAfter $starlPara
$p.Style = "Title"
Add-Run $p "Monday "
Add-Run $p "28/03 02:30pm"

# 2) Heading1 paragraph: Attendees
$prev = $p
$p = New-ParaAfter $prev
$p.Style = "Heading1"
Add-Run $p "Attendees: Calum, Chris, Sergio, Mab, Andreea & Robert"

# 3) Normal paragraph: "Pm meeting:"
$prev = $p
$p = New-ParaAfter $prev
$p.Style = "Normal"
Add-Run $p "Pm meeting:"

# 4) List paragraph: "No timeslot yet for presentation"
$prev = $p
$p = New-ParaAfter $prev
$p.Style = "ListParagraph"
$p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 2, $false, $false)
Add-Run $p "No timeslot yet for presentation"

# 5) List paragraph: "Seminar for " + "STARL next week"
$prev = $p
$p = New-ParaAfter $prev
$p.Style = "ListParagraph"
$p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 2, $false, $false)
Add-Run $p "Seminar for "
Add-Run $p "STARL next week"

# 6) List paragraph: "One more lecture next " + "M" + "on" + "day"
$prev = $p
$p = New-ParaAfter $prev
$p.Style = "ListParagraph"
$p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 2, $false, $false)
Add-Run $p "One more lecture next "
Add-Run $p "M"
Add-Run $p "on"
Add-Run $p "day"

# 7) Normal paragraph: "Group Meeting " + "N" + "otes:"
$prev = $p
$p = New-ParaAfter $prev
$p.Style = "Normal"
Add-Run $p "Group Meeting "
Add-Run $p "N"
Add-Run $p "otes:"

# 8) List paragraph: "Calum " + "& Sergio " + "started implementing questions"
$prev = $p
$p = New-ParaAfter $prev
$p.Style = "ListParagraph"
$p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 2, $false, $false)
Add-Run $p "Calum "
Add-Run $p "& Sergio "
Add-Run $p "started implementing questions"

# 9) List paragraph: "About page content created "
$prev = $p
$p = New-ParaAfter $prev
$p.Style = "ListParagraph"
$p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 2, $false, $false)
Add-Run $p "About page content created "

# 10) List paragraph: "Fixing next buttons "
$prev = $p
$p = New-ParaAfter $prev
$p.Style = "ListParagraph"
$p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 2, $false, $false)
Add-Run $p "Fixing next buttons "

# 11) List paragraph: "Leader board page not updating in real time (need to refresh)"
$prev = $p
$p = New-ParaAfter $prev
$p.Style = "ListParagraph"
$p.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 2, $false, $false)
Add-Run $p "Leader board page not updating in real time (need to refresh)"

# 12) Trailing empty paragraph
$prev = $p
$p = New-ParaAfter $prev
$p.Style = "Normal"
